# Swap the contents of columns B:AD between each pair of adjacent rows.
# (Column A, the running index, stays put on each row.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(33, 34),
    @(158, 159),
    @(181, 182),
    @(300, 301),
    @(350, 351),
    @(368, 369),
    @(371, 372)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    $vals1 = $range1.Value()
    $vals2 = $range2.Value()

    $range1.Value = $vals2
    $range2.Value = $vals1
}
